# Update countries & provincias Spain
# Applies the 12-Apr-2020 22:52 data refresh to the "Pais" sheet:
#  - Updates the "last updated" timestamp string
#  - Refreshes case-count figures for several countries
#  - Re-ranks Peru above Ecuador/Chile, and Niger above Uruguay/Burkina Faso,
#    since their updated totals now place them higher in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados a 12 de Abril de 2020 a las 22:22" -> "...22:52"
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 22:52"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 557071
$ws.Cells.Item(4, 3).Value = 24192
$ws.Cells.Item(4, 5).Value = 503750

# --- Row 17: Brasil ---
$ws.Cells.Item(17, 2).Value = 22169
$ws.Cells.Item(17, 3).Value = 1207
$ws.Cells.Item(17, 5).Value = 20773
$ws.Cells.Item(17, 7).Value = 83
$ws.Cells.Item(17, 8).Value = 1223

# --- Rows 26-28: Peru overtakes Ecuador and Chile in the ranking ---
# Row 26 becomes Peru, with its refreshed totals
$ws.Cells.Item(26, 1).Value = "Peru"
$ws.Cells.Item(26, 2).Value = 7519
$ws.Cells.Item(26, 3).Value = 671
$ws.Cells.Item(26, 4).Value = 1798
$ws.Cells.Item(26, 5).Value = 5528
$ws.Cells.Item(26, 6).Value = 134
$ws.Cells.Item(26, 7).Value = 12
$ws.Cells.Item(26, 8).Value = 193

# Row 27 becomes Ecuador (its totals are unchanged, just shifted down a row)
$ws.Cells.Item(27, 1).Value = "Ecuador"
$ws.Cells.Item(27, 2).Value = 7466
$ws.Cells.Item(27, 3).Value = 209
$ws.Cells.Item(27, 4).Value = 501
$ws.Cells.Item(27, 5).Value = 6632
$ws.Cells.Item(27, 6).Value = 184
$ws.Cells.Item(27, 7).Value = 18
$ws.Cells.Item(27, 8).Value = 333

# Row 28 becomes Chile (its totals are unchanged, just shifted down a row)
$ws.Cells.Item(28, 1).Value = "Chile"
$ws.Cells.Item(28, 2).Value = 7213
$ws.Cells.Item(28, 3).Value = 286
$ws.Cells.Item(28, 4).Value = 2059
$ws.Cells.Item(28, 5).Value = 5074
$ws.Cells.Item(28, 6).Value = 387
$ws.Cells.Item(28, 7).Value = 7
$ws.Cells.Item(28, 8).Value = 80

# --- Row 76: Uzbekistan ---
$ws.Cells.Item(76, 4).Value = 66
$ws.Cells.Item(76, 5).Value = 795

# --- Rows 94-96: Niger overtakes Uruguay and Burkina Faso in the ranking ---
# Row 94 becomes Niger, with its refreshed totals
$ws.Cells.Item(94, 1).Value = "Niger"
$ws.Cells.Item(94, 2).Value = 529
$ws.Cells.Item(94, 3).Value = 38
$ws.Cells.Item(94, 4).Value = 75
$ws.Cells.Item(94, 5).Value = 442
$ws.Cells.Item(94, 6).Value = 0
$ws.Cells.Item(94, 7).Value = 1
$ws.Cells.Item(94, 8).Value = 12

# Row 95 becomes Uruguay (its totals are unchanged, just shifted down a row)
$ws.Cells.Item(95, 1).Value = "Uruguay"
$ws.Cells.Item(95, 2).Value = 501
$ws.Cells.Item(95, 3).Value = 7
$ws.Cells.Item(95, 4).Value = 224
$ws.Cells.Item(95, 5).Value = 270
$ws.Cells.Item(95, 6).Value = 17
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 7

# Row 96 becomes Burkina Faso (its totals are unchanged, just shifted down a row)
$ws.Cells.Item(96, 1).Value = "Burkina Faso"
$ws.Cells.Item(96, 2).Value = 497
$ws.Cells.Item(96, 3).Value = 13
$ws.Cells.Item(96, 4).Value = 161
$ws.Cells.Item(96, 5).Value = 309
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 27

# --- Row 106: Nigeria ---
$ws.Cells.Item(106, 2).Value = 323
$ws.Cells.Item(106, 3).Value = 5
$ws.Cells.Item(106, 4).Value = 85
$ws.Cells.Item(106, 5).Value = 228

# --- Row 110: Montenegro ---
$ws.Cells.Item(110, 2).Value = 272
$ws.Cells.Item(110, 3).Value = 9
$ws.Cells.Item(110, 5).Value = 264
